$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 2; D = '67.037.65'; E = '  +0.51%  ' },
    @{ Row = 3; D = '3.472.81'; E = '  +0.82%  ' },
    @{ Row = 4; E = '  -0.05%  ' },
    @{ Row = 5; D = '587.23'; E = '  -0.55%  ' },
    @{ Row = 6; D = '177.51'; E = '  -0.05%  ' },
    @{ Row = 7; E = '  -0.03%  ' },
    @{ Row = 8; D = '0.601'; E = '  -1.95%  ' },
    @{ Row = 9; D = '3.468.31'; E = '  +0.57%  ' },
    @{ Row = 10; D = '0.133'; E = '  -2.50%  ' },
    @{ Row = 11; D = '6.91'; E = '  -0.21%  ' },
    @{ Row = 12; D = '0.422'; E = '  -1.91%  ' },
    @{ Row = 13; D = '4.078.76'; E = '  +0.68%  ' },
    @{ Row = 14; D = '30.60'; E = '  -2.95%  ' },
    @{ Row = 15; E = '  -0.57%  ' },
    @{ Row = 16; D = '66.948.73'; E = '  +0.24%  ' },
    @{ Row = 17; D = '0.0000173'; E = '  -1.08%  ' },
    @{ Row = 18; D = '3.481.90'; E = '  +0.89%  ' },
    @{ Row = 19; D = '6.00'; E = '  -3.51%  ' },
    @{ Row = 20; D = '13.92'; E = '  -1.13%  ' },
    @{ Row = 21; D = '382.02'; E = '  -1.51%  ' },
    @{ Row = 22; D = '7.86'; E = '  -0.27%  ' },
    @{ Row = 23; B = 'Litecoin'; C = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D = '72.66'; E = '  +1.19%  ' },
    @{ Row = 24; B = 'Polygon'; C = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D = '0.541'; E = '  +1.29%  ' },
    @{ Row = 25; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '0.999'; E = '  +0.05%  ' },
    @{ Row = 26; B = 'LEO'; C = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; D = '5.74'; E = '  -0.27%  ' },
    @{ Row = 27; E = '  +0.72%  ' },
    @{ Row = 28; D = '9.89'; E = '  -3.18%  ' },
    @{ Row = 29; E = '  +1.83%  ' },
    @{ Row = 30; E = '  +0.55%  ' },
    @{ Row = 31; D = '24.35'; E = '  +4.41%  ' },
    @{ Row = 32; D = '5.89'; E = '  -4.10%  ' },
    @{ Row = 33; E = '  -1.83%  ' },
    @{ Row = 34; D = '1.33'; E = '  -4.01%  ' },
    @{ Row = 35; D = '1.00'; E = '  +0.08%  ' },
    @{ Row = 36; D = '7.20'; E = '  -1.42%  ' },
    @{ Row = 37; E = '  +0.99%  ' },
    @{ Row = 38; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '160.84'; E = '  -1.32%  ' },
    @{ Row = 39; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '29.28'; E = '  +12.26%  ' },
    @{ Row = 40; D = '0.894'; E = '  +2.46%  ' },
    @{ Row = 41; B = 'dogwifhat'; C = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D = '2.66'; E = '  -3.43%  ' },
    @{ Row = 42; B = 'Stacks'; C = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D = '1.80'; E = '  -3.04%  ' },
    @{ Row = 43; D = '4.52'; E = '  -2.08%  ' },
    @{ Row = 44; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D = '6.55'; E = '  -3.47%  ' },
    @{ Row = 45; B = 'Maker'; C = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D = '2.764.50'; E = '  +1.25%  ' },
    @{ Row = 46; D = '0.0699'; E = '  -2.52%  ' },
    @{ Row = 47; D = '40.84'; E = '  -0.43%  ' },
    @{ Row = 48; D = '24.80'; E = '  -4.66%  ' },
    @{ Row = 49; D = '0.0296'; E = '  -0.36%  ' },
    @{ Row = 50; D = '311.46'; E = '  -3.52%  ' },
    @{ Row = 51; D = '1.01'; E = '  -2.93%  ' }
)

foreach ($change in $changes) {
    $r = $change.Row
    if ($change.ContainsKey('B')) {
        $ws.Cells.Item($r, 2).Value = $change.B
    }
    if ($change.ContainsKey('C')) {
        $ws.Cells.Item($r, 3).Value = $change.C
    }
    if ($change.ContainsKey('D')) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $change.D
    }
    if ($change.ContainsKey('E')) {
        $ws.Cells.Item($r, 5).Value = $change.E
    }
}
